$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.372.93'
$ws.Range('E2').Value = '  -3.03%  '
$ws.Range('D3').Value = '2.978.36'
$ws.Range('E3').Value = '  -4.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '495.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '2.975.94'
$ws.Range('E8').Value = '  -4.93%  '
$ws.Range('E9').Value = '  -3.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.28'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('E11').Value = '  -3.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.352'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.40%  '
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '3.484.64'
$ws.Range('E14').Value = '  -4.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.86'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.28%  '
$ws.Range('D16').Value = '56.256.39'
$ws.Range('E16').Value = '  -3.20%  '
$ws.Range('D17').Value = '2.973.89'
$ws.Range('E17').Value = '  -4.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000146'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.81'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.72'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '324.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.33%  '
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.465'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '61.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -9.82%  '
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.161'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.06%  '
$ws.Range('E28').Value = '  -5.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.47'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.72'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('E32').Value = '  -2.95%  '
$ws.Range('E33').Value = '  -7.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.03'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '154.87'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.41%  '
$ws.Range('E37').Value = '  -6.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.58'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -9.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0677'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.29'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.60%  '
$ws.Range('D41').Value = '3.005.24'
$ws.Range('E41').Value = '  -4.94%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.42'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -9.62%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.997'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.636'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.71%  '
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.994'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.47%  '
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('D47').Value = '2.207.45'
$ws.Range('E47').Value = '  -2.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.55'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -8.45%  '
$ws.Range('E49').Value = '  +5.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0237'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.68%  '
